# Apply updated dSF (column F) values per repull / mean calculation update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    3  = -6
    5  = -4
    6  = -1
    9  = -3
    16 = -6
    17 = -5
    18 = 5
    19 = 2
    20 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
